$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the C1 header: "pob _libre" -> "pob_libre"
$ws.Range("C1").Value = "pob_libre"

# Move the active selection to C1 (was B13)
$ws.Range("C1").Select()
